$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.934.98"
$ws.Range("E2").Value = "  -2.36%  "

# Row 3
$ws.Range("D3").Value = "3.578.15"
$ws.Range("E3").Value = "  -3.48%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.78%  "

# Row 7
$ws.Range("D7").Value = "3.574.42"
$ws.Range("E7").Value = "  -3.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.97%  "

# Row 9
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.676"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.90%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.79%  "

# Row 15
$ws.Range("D15").Value = "4.151.14"
$ws.Range("E15").Value = "  -3.34%  "

# Row 16
$ws.Range("D16").Value = "3.582.62"
$ws.Range("E16").Value = "  -3.28%  "

# Row 17
$ws.Range("E17").Value = "  -1.46%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.86%  "

# Row 19
$ws.Range("D19").Value = "66.938.13"
$ws.Range("E19").Value = "  -2.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.85%  "

# Row 21
$ws.Range("E21").Value = "  -7.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.84%  "

# Row 28
$ws.Range("E28").Value = "  +0.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.22%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.49%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "629.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.80%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.73%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.399"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.38%  "

# Row 39
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -7.52%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.148.29"
$ws.Range("E41").Value = "  +7.57%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.133"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.30%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0414"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.92%  "

# Row 47
$ws.Range("E47").Value = "  +1.49%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.91%  "
